$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the daily series: insert a new row
# at position 63 (this shifts the former rows 63..148 down to 64..149,
# matching every row below it).
$ws.Rows.Item(63).Insert()

# Fill the newly inserted row 63 with its data. It reuses the same
# market/product/quality/origin/unit metadata as the row that used to
# sit at position 63 (now at 64), but carries its own date, volume,
# prices and $/kg value.
$ws.Cells.Item(63, 1).Value = 6
$ws.Cells.Item(63, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(63, 3).Value = "Metropolitana"
$ws.Cells.Item(63, 4).Value = 44579
$ws.Cells.Item(63, 5).Value = 13
$ws.Cells.Item(63, 6).Value = "Fruta"
$ws.Cells.Item(63, 7).Value = 100101
$ws.Cells.Item(63, 8).Value = "Berries"
$ws.Cells.Item(63, 9).Value = 100101004
$ws.Cells.Item(63, 10).Value = "Frambuesa"
$ws.Cells.Item(63, 11).Value = "Sin especificar"
$ws.Cells.Item(63, 12).Value = "Especial"
$ws.Cells.Item(63, 13).Value = 250
$ws.Cells.Item(63, 14).Value = 8000
$ws.Cells.Item(63, 15).Value = 8000
$ws.Cells.Item(63, 16).Value = 8000
$ws.Cells.Item(63, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(63, 18).Value = "Provincia de Linares"
$ws.Cells.Item(63, 19).Value = 4000
$ws.Cells.Item(63, 20).Value = 2

# The D column (Fecha) carries a date-aware number format in this
# sheet; apply the same format used by every other row's date cell.
$ws.Cells.Item(63, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
